# MCP33151 Eval Board Schematic Symbol and Footprint Checks
# - Add a "Check before order" sheet (copy of Sheet1) with updated component
#   names/notes and a few newly-checked parts appended.
# - Sheet1's selection moves off E13 (no longer the active tab) to A11.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Park Sheet1's selection where it ends up in the final file before we move
# focus to the new sheet (keeps Sheet1 from being "tabSelected").
$ws1.Range("A11").Select()

# Duplicate Sheet1 (keeps formatting, merged cells and the picture) and place
# the copy right after it, then rename it.
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Check before order"

# Widen column A a bit to fit the longer component descriptions.
$ws2.Columns.Item(1).ColumnWidth = 37

# Row 2: TPS65131 -> "TPS65131 boost"; clear the old Yes-markers in B:E.
$ws2.Range("A2").Value = "TPS65131 boost"
$ws2.Range("B2:E2").ClearContents()

# Row 3: MC79L12ACDR -> "MC79L12ACDR -12v"; clear B3, keep the merged
# Altium Vault note in C3:E7.
$ws2.Range("A3").Value = "MC79L12ACDR -12v"
$ws2.Range("B3").ClearContents()

# Row 4: TLV76012DBZR -> "TLV76012DBZR +12v"; clear B4.
$ws2.Range("A4").Value = "TLV76012DBZR +12v"
$ws2.Range("B4").ClearContents()

# Row 5: TLV74118PDBVR -> "TLV74118PDBVR +1v8"; clear B5.
$ws2.Range("A5").Value = "TLV74118PDBVR +1v8"
$ws2.Range("B5").ClearContents()

# Row 6: ADR5045BRTZ-REEL7 -> "ADR5045BRTZ-REEL7 +5v_a"; clear B6.
$ws2.Range("A6").Value = "ADR5045BRTZ-REEL7 +5v_a"
$ws2.Range("B6").ClearContents()

# Row 7: MCP9700T-E/TT -> "MCP9700T-E/TT temp sensor"; clear B7.
$ws2.Range("A7").Value = "MCP9700T-E/TT temp sensor"
$ws2.Range("B7").ClearContents()

# Rows 8-13 keep their component names, but the per-column Yes checkmarks
# (B:E) haven't been re-verified yet on this checklist, so clear them -
# except row 13's "Yes" in column B, which stays.
$ws2.Range("B8:E8").ClearContents()
$ws2.Range("B9:E9").ClearContents()
$ws2.Range("B10:E10").ClearContents()
$ws2.Range("B11:E11").ClearContents()
$ws2.Range("B12:E12").ClearContents()
$ws2.Range("C13:E13").ClearContents()

# New parts added to the board - append them below the existing checklist.
$ws2.Range("A14").Value = "1985195 terminal blocks"
$ws2.Range("A15").Value = "MCP23008-E/ML io expansion"
$ws2.Range("A16").Value = "BAT54s"
$ws2.Range("A17").Value = "INA821ID in-amp"
$ws2.Range("A18").Value = "AO3407A"
$ws2.Range("B18").Value = "Yes"
$ws2.Range("C18").Value = "Yes"
$ws2.Range("D18").Value = "Yes"
$ws2.Range("E18").Value = "Yes"

# Make the new sheet active and leave the selection where the author left it.
$ws2.Activate()
$ws2.Range("B27").Select()
